# Updates for Framework and CVC Project up to 6/19/2018
#
# Before: messageSchema, sample01 (old data), sample02 (old data)
# After:  messageSchema, sample01 (new geolocation schema), sampleold (= old sample01 data)
# (old sample02 data is dropped entirely)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Re-arrange sheets.
#    - old "sample01" keeps its data/format, just gets renamed to "sampleold"
#    - old "sample02" is discarded completely
#    - a brand-new "sample01" sheet is inserted before "sampleold" with the
#      new geolocation-oriented schema
# ---------------------------------------------------------------------------
$oldSample01 = $wb.Worksheets.Item("sample01")
$oldSample01.Name = "sampleold"

$oldSample02 = $wb.Worksheets.Item("sample02")
$oldSample02.Delete()

$newSample01 = $wb.Worksheets.Add($wb.Worksheets.Item("sampleold"))
$newSample01.Name = "sample01"

# ---------------------------------------------------------------------------
# 2. Populate the new "sample01" sheet
# ---------------------------------------------------------------------------
$newSample01.Columns.Item(1).ColumnWidth = 24.3046875
$newSample01.Columns.Item(2).ColumnWidth = 19
$newSample01.Columns.Item(3).ColumnWidth = 19
$newSample01.Columns.Item(4).ColumnWidth = 14.3046875
$newSample01.Columns.Item(5).ColumnWidth = 20.15234375
$newSample01.Columns.Item(6).ColumnWidth = 13.53515625
$newSample01.Columns.Item(7).ColumnWidth = 12.53515625
$newSample01.Columns.Item(8).ColumnWidth = 12.53515625
$newSample01.Columns.Item(9).ColumnWidth = 19.23046875

$newSample01.Range("A1").Value = "FieldName"
$newSample01.Range("B1").Value = "timestamp"
$newSample01.Range("C1").Value = "targetSite"
$newSample01.Range("D1").Value = "deviceId"
$newSample01.Range("E1").Value = "eventType"
$newSample01.Range("F1").Value = "sourceName"
$newSample01.Range("G1").Value = "latitude"
$newSample01.Range("H1").Value = "longitude"
$newSample01.Range("I1").Value = "activity"

$newSample01.Range("A2").Value = "FieldType"
$newSample01.Range("B2").Value = "timestamp"
$newSample01.Range("C2").Value = "string"
$newSample01.Range("D2").Value = "string"
$newSample01.Range("E2").Value = "string"
$newSample01.Range("F2").Value = "string"
$newSample01.Range("G2").Value = "string"
$newSample01.Range("H2").Value = "string"
$newSample01.Range("I2").Value = "integer"

$newSample01.Range("A3").Value = "FieldUsage"
$newSample01.Range("B3").Value = "auto"
$newSample01.Range("C3").Value = "key"
$newSample01.Range("D3").Value = "key"
$newSample01.Range("E3").Value = "key"

$newSample01.Range("A4").Value = "absoluteMinValue"

$newSample01.Range("A5").Value = "absoluteMaxValue"

$newSample01.Range("A6").Value = "normalMinValue"

$newSample01.Range("A7").Value = "normalMaxValue"

$newSample01.Range("A8").Value = "message"
$newSample01.Range("C8").Value = "cvcdev"
$newSample01.Range("D8").Value = "locdevsrc"
$newSample01.Range("E8").Value = "geolocation"
$newSample01.Range("F8").Value = "Prime Gas"
$newSample01.Range("G8").Value = 77.3451
$newSample01.Range("H8").Value = 80.3456
$newSample01.Range("I8").Value = 1

$newSample01.Range("A1:I8").Select()
$newSample01.Range("E12").Select()

# Freeze header row + first column like the other data sheets
$newSample01.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 3. Trim the "messageSchema" sheet down to the new (fewer) columns
# ---------------------------------------------------------------------------
$schema = $wb.Worksheets.Item("messageSchema")

$schema.Range("E1").Value = "eventType"
$schema.Range("F1").Value = "sourceName"
$schema.Range("G1").Value = "latitude"
$schema.Range("H1").Value = "longitude"
$schema.Range("I1").Value = "activity"
$schema.Range("J1:L1").ClearContents()

$schema.Range("G2").Value = "string"
$schema.Range("H2").Value = "string"
$schema.Range("I2").Value = "string"
$schema.Range("J2:L2").ClearContents()

$schema.Range("E3").Value = "key"

$schema.Range("G4:L4").ClearContents()
$schema.Range("G5:L5").ClearContents()
$schema.Range("G6:L6").ClearContents()
$schema.Range("G7:L7").ClearContents()

$schema.Activate()
$schema.Range("B2").Select()
$excel.ActiveWindow.FreezePanes = $true

# ---------------------------------------------------------------------------
# 4. Leave "sample01" (the new sheet) as the active/selected tab, matching
#    the workbook's activeTab pointing at the second sheet.
# ---------------------------------------------------------------------------
$newSample01.Activate()
$newSample01.Range("E12").Select()
